# Fill in the missing "1" flags across the pathways sheet and refresh the
# saved view/selection state, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where columns D, E, F need a value of 1 (column C already holds a 1
# for these rows).
$rowsDEF = @(2,3,4,5,6,7,8,9,10,11,12,13,16,17,18,19,20,21,22,23,24,25,26,31,44)

foreach ($r in $rowsDEF) {
    $ws.Cells.Item($r, 4).Value = 1   # D
    $ws.Cells.Item($r, 5).Value = 1   # E
    $ws.Cells.Item($r, 6).Value = 1   # F
}

# Rows where columns C, D, E, F all need a value of 1 (C was previously
# blank on these rows too).
$rowsCDEF = @(27,30,32,45,46,59,72,73)

foreach ($r in $rowsCDEF) {
    $ws.Cells.Item($r, 3).Value = 1   # C
    $ws.Cells.Item($r, 4).Value = 1   # D
    $ws.Cells.Item($r, 5).Value = 1   # E
    $ws.Cells.Item($r, 6).Value = 1   # F
}

# Refresh the window/view state: move the selection to F91, then scroll so
# row 57 becomes the top visible row.
$ws.Range("F91").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 57
$win.ScrollColumn = 1
